$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gas_criteria")

# Row 2: ABG threshold positive
$ws.Range("B2").Value = 526
$ws.Range("C2").Value = 26.5

# Row 3: VBG threshold positive
$ws.Range("B3").Value = 865
$ws.Range("C3").Value = 43.6

# Row 4: PCO2 OTHER -> PCO2 UNKNOWN threshold positive
$ws.Range("A4").Value = "PCO2 UNKNOWN threshold positive"
$ws.Range("B4").Value = 69
$ws.Range("C4").Value = 3.5

# Row 5: Any gas threshold positive
$ws.Range("B5").Value = 1077
$ws.Range("C5").Value = 54.3
